$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$rsquo = [char]0x2019

# Update top players ranking requirement text (row 15): "top ten" -> "top five"
$ws.Range("A15").Value = "The main page MUST display the top five players and the current player" + $rsquo + "s ranking."

# Update avatar-name correction requirement text (row 8)
$ws.Range("A8").Value = "Users MUST be able to correct avatar names."

# Update the visible window/selection state to match the saved view
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("A9").Select()
